$d = $word.ActiveDocument

$d.Content.Find.Execute("{Data raskida}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{termDate}", 2)
